$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update CR_REVIEW_001 (row 2) and CR_REVIEW_003 (row 4) Status to "Closed"
$ws.Range("H2").Value = "Closed"
$ws.Range("H4").Value = "Closed"

# Reflect the scrolled / selected view state recorded in the saved file
try {
    $excel.ActiveWindow.ScrollColumn = 6
} catch {
}
$ws.Range("H3").Select()
